# tien commit sua thong bao
# Update the "Mã số" (ID) column: 0017-0020 -> 0021-0024

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "0021"
$ws.Range("A3").Value = "0022"
$ws.Range("A4").Value = "0023"
$ws.Range("A5").Value = "0024"

# Move the active selection to E7 (as in the saved workbook state)
$ws.Range("E7").Select()
